$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.831.02"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "1.895.74"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3174"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07052"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08061"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7729"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.50%  "

$ws.Range("D13").Value = "1.889.86"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.348"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "29.831.63"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.033"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007720"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.289"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.18%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "2.140.68"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("E24").Value = "  -0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1662"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.354"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.88%  "

$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.062"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.408"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.539"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.442"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05727"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.052"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.262"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("E36").Value = "  +1.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9985"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.633"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01913"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.786"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.819"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8428"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").Value = "1.035.10"
$ws.Range("E46").Value = "  +4.57%  "

$ws.Range("E47").Value = "  +2.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.875"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.978"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.442"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.71%  "

$ws.Range("D51").Value = "2.031.34"
$ws.Range("E51").Value = "  -0.43%  "
